# Build site at 2023-04-12 14:53:07 UTC
#
# LOB1024.xlsx: the "Objetivos/Docentes/Programa/Avaliacao" block had its
# B/C content column out of sync with the A label column (stale / duplicated
# placeholder values). This fixes the data by:
#   1. Inserting one new row at row 13 to hold the "Docentes responsaveis:"
#      value (which previously had no row of its own), pushing every row
#      below it down by one.
#   2. Filling in the correct B/C content for the rows whose value was
#      wrong/missing, reusing the already-correct content where present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at 13 - shifts old rows 13..24 down to 14..25,
#    carrying their values AND row heights with them.
$ws.Rows.Item(13).Insert()

# The row-13 insert leaves a blank styled "A13" cell behind (copied from the
# row above); the target layout has no A-cell at all on this row, so drop it
# completely (Clear removes the cell, not just its contents).
$ws.Range("A13").Clear()

# 2) Objetivos: (row 10) previously held the wrong "Docentes" text - replace
#    with the real objectives paragraph.
$ws.Range("B10").Value = "Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na Estática dos Sólidos"
$ws.Range("C10").Value = "Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na Estática dos Sólidos"

# 3) The newly inserted row 13 is the "Docentes responsáveis:" value (no
#    label in column A). Pick up the normal B/C formatting (regular black
#    text + wrap for B, red "modified" text + wrap for C) from row 10 before
#    filling in the value, since a freshly-inserted row has no B/C cells yet.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840650 - Janaína Ferreira Batista"
$ws.Range("C13").Value = "5840650 - Janaína Ferreira Batista"

# 4) Programa resumido: (row 14 after the shift) previously held the
#    placeholder "Semestral" - replace with the real short syllabus (PT).
$ws.Range("B14").Value = "Estática das Partículas, Estática de Corpos Rígidos, Equilíbrio de Corpos Rígidos, Análise de Estruturas, Forças Distribuídas, Forças em Vigas."
$ws.Range("C14").Value = "Estática das Partículas, Estática de Corpos Rígidos, Equilíbrio de Corpos Rígidos, Análise de Estruturas, Forças Distribuídas, Forças em Vigas."

# 5) Programa: (row 16 after the shift) previously held a stray duplicate
#    "01/01/2018" - replace with the full syllabus text (PT).
$ws.Range("B16").Value = "Programa: 1) Estática de partículas: Vetores. Resultante de várias forças concorrentes. Equilíbrio de uma partícula. 2) Estática de Corpos Rígidos: Conceito de corpo rígido. Forças externas e forças internas. Forças equivalentes. Momento de uma força com relação a um ponto. Sistemas equivalentes de forças. Diagrama de corpo livre.3) Equilíbrio de corpos rígidos: Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações de apoios e conexões para uma estrutura 3D. Equilíbrio de um corpo rígido em 3D. 4) Análise de Estruturas: Treliças: Definições. Treliça simples. Análise de treliças pelo método dos nós. Análise de treliças pelo método das seções. Estruturas: estruturas que contêm elementos sujeitos a ação de múltiplas forças, transmissão e modificação de forças.5) Forças Distribuídas: Centróides e baricentros de áreas, linhas e volumes; teoremas de Guldinus-Pappus; cargas distribuídas sobre vigas.6) Forças em Vigas: Diagramas de cisalhamento e momento fletor para uma viga carregada."
$ws.Range("C16").Value = "Programa: 1) Estática de partículas: Vetores. Resultante de várias forças concorrentes. Equilíbrio de uma partícula. 2) Estática de Corpos Rígidos: Conceito de corpo rígido. Forças externas e forças internas. Forças equivalentes. Momento de uma força com relação a um ponto. Sistemas equivalentes de forças. Diagrama de corpo livre.3) Equilíbrio de corpos rígidos: Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações de apoios e conexões para uma estrutura 3D. Equilíbrio de um corpo rígido em 3D. 4) Análise de Estruturas: Treliças: Definições. Treliça simples. Análise de treliças pelo método dos nós. Análise de treliças pelo método das seções. Estruturas: estruturas que contêm elementos sujeitos a ação de múltiplas forças, transmissão e modificação de forças.5) Forças Distribuídas: Centróides e baricentros de áreas, linhas e volumes; teoremas de Guldinus-Pappus; cargas distribuídas sobre vigas.6) Forças em Vigas: Diagramas de cisalhamento e momento fletor para uma viga carregada."

# 6) Método: (row 19 after the shift) previously held the wrong "Docentes"
#    text again - replace with the real grading-method paragraph (this was
#    already-correct text that lived one row down, attached to the wrong
#    label, before the fix).
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# 7) Critério: (row 20 after the shift) - move up the "NF>=5,0." text that
#    had been sitting one row too low.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# 8) Norma de recuperação: (row 21 after the shift) - move up the
#    "(NF+RC)/2..." text that had been sitting one row too low.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# 9) Bibliografia: (row 22 after the shift) was empty - fill in the
#    bibliography text.
$ws.Range("B22").Value = "1. BEER, Ferdinand Pierre, ; JOHNSTON, E. Russel.; Eisenberg, Elliot R., Mecânica vetorial para engenheiros: Estática.  Mc Graw Hill (2011).2. HIBBELER, R.C. Mecânica para engenharia, Vol. 1: estática, Pearson Prentice Hall (2005).3. MERIAM J. L. ; KRAIGE, L. G., Mecânica, estática, Livros Técnicos e Científicos Editora (2004)."
$ws.Range("C22").Value = "1. BEER, Ferdinand Pierre, ; JOHNSTON, E. Russel.; Eisenberg, Elliot R., Mecânica vetorial para engenheiros: Estática.  Mc Graw Hill (2011).2. HIBBELER, R.C. Mecânica para engenharia, Vol. 1: estática, Pearson Prentice Hall (2005).3. MERIAM J. L. ; KRAIGE, L. G., Mecânica, estática, Livros Técnicos e Científicos Editora (2004)."

# 10) Column layout cleanup: column A's width declaration no longer spans
#     into column B (columns A and B have been independently-declared
#     widths all along; this just tightens the metadata so col A stands on
#     its own single-column range).
$ws.Columns.Item(1).ColumnWidth = 30.7109375
